# "Development Center in Extension"
# Adds a new top-level DEV-MENU entry ("build modeling / development
# center") with four child entries (entity / model / relation / index),
# and demotes the previous "modeling management" entry to be a child of
# the new top-level entry (relabelled as the legacy/old entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Title bar (C2:J2) gets a bigger, red font (new font + style).
# ---------------------------------------------------------------------
$ws.Range("C2:J2").Font.Size = 16
$ws.Range("C2:J2").Font.Color = 255

# ---------------------------------------------------------------------
# 2. Grow the table: clone row 5's layout/format down into rows 6-10
#    so every new row starts from the same base formatting.
# ---------------------------------------------------------------------
$ws.Range("A5:J5").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A9:J9").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A10:J10").PasteSpecial(-4122)

# Row 6 is the new top-level entry: highlight its order/level cells and
# give its uri cell the bold "EXPAND" look (new red font, no fill).
$ws.Range("A4").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Font.Size = 16
$ws.Range("I6").Font.Color = 255

# ---------------------------------------------------------------------
# 3. Row 6 - new top-level "DEV-MENU" entry.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "a344bb6e-a669-4f99-93ad-7f73ecd5bc6b"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "DEV-MENU"
$ws.Range("D6").Value = 20000
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "develop.atom"
$ws.Range("G6").Value = "建模管理"
$ws.Range("H6").Value = "apartment"
$ws.Range("I6").Value = "EXPAND"

# ---------------------------------------------------------------------
# 4. Row 5 - existing entry, re-parented under the new top-level entry
#    and relabelled as the legacy/old entry.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "ee0a6f15-62cf-42fd-b6dd-c1043e232486"
$ws.Range("B5").Formula = "=A$6"
$ws.Range("C5").Value = "DEV-MENU"
$ws.Range("D5").Value = 2000
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "zero.develop.model"
$ws.Range("G5").Value = "「旧」建模管理"
$ws.Range("H5").Value = "appstore"
$ws.Range("I5").Value = "/epic/modeling"

# ---------------------------------------------------------------------
# 5. Rows 7-10 - new child entries under the new top-level entry.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "3f6d8b63-f569-4b5c-b109-48e74b3661ce"
$ws.Range("B7").Formula = "=A$6"
$ws.Range("C7").Value = "DEV-MENU"
$ws.Range("D7").Value = 1005
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "develop.atom.entity"
$ws.Range("G7").Value = "表实体"
$ws.Range("H7").Value = "table"
$ws.Range("I7").Value = "/atom/entity"

$ws.Range("A8").Value = "c5389bd0-a37e-4c6b-b0c2-c499940f498a"
$ws.Range("B8").Formula = "=A$6"
$ws.Range("C8").Value = "DEV-MENU"
$ws.Range("D8").Value = 1010
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "develop.atom.model"
$ws.Range("G8").Value = "领域模型"
$ws.Range("H8").Value = "experiment"
$ws.Range("I8").Value = "/atom/model"

$ws.Range("A9").Value = "ad1075e4-4bf5-4ef2-a48f-2e1e5c60305d"
$ws.Range("B9").Formula = "=A$6"
$ws.Range("C9").Value = "DEV-MENU"
$ws.Range("D9").Value = 1015
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "develop.atom.relation"
$ws.Range("G9").Value = "关系设置"
$ws.Range("H9").Value = "share-alt"
$ws.Range("I9").Value = "/atom/relation"

$ws.Range("A10").Value = "1d0cda23-e8db-40b1-abe0-75f748b04478"
$ws.Range("B10").Formula = "=A$6"
$ws.Range("C10").Value = "DEV-MENU"
$ws.Range("D10").Value = 1020
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "develop.atom.index"
$ws.Range("G10").Value = "索引设置"
$ws.Range("H10").Value = "search"
$ws.Range("I10").Value = "/atom/indexing"

# ---------------------------------------------------------------------
# 6. Selection / view bookkeeping to match the edited workbook's state.
# ---------------------------------------------------------------------
$ws.Range("D7").Select()

Write-Output "done"
